# PO2EBL_SIQ.xlsx edit — "Moving SIQ to solve folder structure violation"
#
# The underlying content change in this commit is:
#   - Header D1 "Question"      -> "Question\n"      (trailing newline added)
#   - Header H1 "Response date" -> "Response date\n" (trailing newline added)
#   - New "Status" values ("accepted ") filled in for the four SIQ rows (J2:J5)
#   - Active window view nudged (zoom back to 100%, selection moved to E5)
#
# (Everything else in the raw XML diff — defaultThemeVersion, absPath,
#  calcId/iterateDelta, theme minor/major CJK font substitutions, exact
#  row-height / column-width fractions, styles.xml namespace additions —
#  is Excel-build/version round-trip noise that isn't driven by any
#  deliberate user action in the Excel object model, so it isn't
#  reproduced here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("D1").Value = "Question`n"
$ws.Range("H1").Value = "Response date`n"

# --- Status column (added "accepted " answers for the 4 SIQ rows) ------
$ws.Range("J2").Value = "accepted "
$ws.Range("J3").Value = "accepted "
$ws.Range("J4").Value = "accepted "
$ws.Range("J5").Value = "accepted "

# --- View state: rezoom to 100% and move the selection to E5 -----------
$excel.ActiveWindow.Zoom = 100
$ws.Range("E5").Select()
